# Append new weekly workout records to the Kilimanjaro scoreboard.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns: A Participant, B Date, C Workout Type, D Total Duration,
#          E Total Distance, F Total Elevation, G Zone1, H Zone2,
#          I Zone3, J Zone4, K Zone5, L Workout Level, M Week

$rows = @(
    @("Steven", 45461, "Walk", 21, 1,                   20,  21, 0,  0,  0,  0, "Agile Antelope", 2),
    @("Matt",   45462, "Run",  20, 2.1800000000000002,  194, 0,  18, 0,  0,  0, "Agile Antelope", 2),
    @("Matt",   45462, "Ride", 32, 0,                   0,   29, 2,  0,  0,  0, "Agile Antelope", 2),
    @("Eric",   45462, "Run",  47, 4.4800000000000004,  151, 0,  2,  15, 26, 0, "Brave Leopard",  2),
    @("Steven", 45462, "Walk", 55, 2.15,                315, 55, 0,  0,  0,  0, "Agile Antelope", 2),
    @("Steven", 45462, "Walk", 48, 1.99,                171, 48, 0,  0,  0,  0, "Agile Antelope", 2),
    @("Matt",   45462, "Walk", 14, 0.55000000000000004, 39,  14, 0,  0,  0,  0, "Agile Antelope", 2),
    @("Eric",   45462, "Ride", 57, 0,                   0,   13, 33, 8,  2,  0, "Brave Leopard",  2)
)

$startRow = 64
$templateRow = 63

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Copy the last existing data row as a template so the new row
    # inherits identical cell styling (e.g. the date-formatted style).
    $ws.Range("A" + $templateRow + ":M" + $templateRow).Copy($ws.Range("A" + $r + ":M" + $r))

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
    $ws.Cells.Item($r, 10).Value = $data[9]
    $ws.Cells.Item($r, 11).Value = $data[10]
    $ws.Cells.Item($r, 12).Value = $data[11]
    $ws.Cells.Item($r, 13).Value = $data[12]
}

# Update the active selection to reflect the new bottom of the data range
# (mirrors the author scrolling down to the newly added rows).
$ws.Range("M67").Select()
